# FINAL DEL 18 SEPT 2021
# Moves the payroll sheet from "SEMANA 37" (06-12 Sep 2021) to
# "SEMANA 38" (13-19 Sep 2021): updates the week-label string, a handful
# of pay-table numbers, and the scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Week label (B9 is the single source; H9/B27/H27/B43/H43/B60 all
#     reference it via formulas, e.g. H9 = "=B9", B27 = "=B9", etc., so
#     they update themselves on recalculation) -------------------------
$ws.Range("B9").Value = "SEMANA   38  DEL  13      Al    19 DE SEPTIEMBRE          2021"

# --- First pay block (rows 3-6): add the 867 "EXTRA" adjustment -------
# K6 (=SUM(K3:K5)) recalculates automatically: 2600 + 867 = 3467
$ws.Range("K4").Value = 867

# --- Second pay block (rows 21-24): drop the hour count, replace the
#     computed "280*J21" total with a flat figure -----------------------
$ws.Range("J21").ClearContents()
$ws.Range("K21").Value = 2800
# K24 (=SUM(K21:K23)) recalculates automatically: 2800 + 0 + 0 = 2800

# --- Third pay block (rows 38-41): zero out the loan/advance line -----
$ws.Range("K40").Value = 0
# K41 (=SUM(K38:K40)) recalculates automatically: 2500 - 500 + 0 = 2000

# --- Scroll the sheet so row 37 is at the top, keeping K41 selected,
#     matching the saved view (topLeftCell changes from A22 to A37) ---
[void]$ws.Range("K41").Select()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 37
$aw.ScrollColumn = 1

# NOTE: the TODAY()/NOW()-derived date cells (C14, I14, C32, I32, C48,
# I48, C65) are left untouched on purpose -- they keep their original
# formulas (=TODAY(), =C14, =I14, ...) and simply recalculate against
# the host clock when the workbook is saved.
